# Reverse the order of comma-separated "Recorded By" entries in column G.
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# Cells holding a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $text = $cell.Text

    if ($text -ne $null -and $text -ne "") {
        $asString = [string]$text
        if ($asString.Contains(",")) {
            $parts = $asString.Split(",")
            $n = $parts.Length
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i].Trim()
            }
            $joined = $reversed -join ", "
            $cell.Value = $joined
        }
    }
}
